# Commit: "added a sentence muhahahha9"
#
# The only meaningful content change is a new sentence inserted (as a
# tracked insertion by "Loni") right after the third "Dumb:" + quote
# occurrence in the document -- the short paragraph that sits right before
# the trailing empty paragraph and the final "Pink sticky change..."
# paragraph. The new text reads:
#     " I choose my CUCUMBER MUHAHA" + closing quote
# so the paragraph ends up saying: Dumb:" I choose my CUCUMBER MUHAHA"
#
# The "_GoBack" bookmark (marking the most recent edit location) also
# moves from the earlier "Fine, choose your weapon!" line to sit inside
# this newly typed sentence, which is what Word does automatically when
# you type/edit at a new spot.

$d = $word.ActiveDocument

# Keep the edit tracked as a revision, same as the rest of the document.
$d.TrackRevisions = $true

# The document re-uses the same left double quotation mark (U+201C) as
# both an opening and a closing quote mark throughout, so build literals
# with that exact character.
$quote = [char]0x201C
$target = "Dumb:" + $quote

# There are three occurrences of "Dumb:" + quote in the body; we need the
# third one.
$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$matchCount = 0
$keepGoing = $true
while ($keepGoing -and $matchCount -lt 3) {
    $keepGoing = $rng.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($keepGoing) {
        $matchCount = $matchCount + 1
        if ($matchCount -lt 3) {
            $rng.Collapse(0)
            $rng.End = $d.Content.End
        }
    }
}

if ($matchCount -eq 3) {
    # Collapse to right after the found "Dumb:"" text and insert the new
    # sentence there, authored by Loni.
    $rng.Collapse(0)
    $insertStart = $rng.Start

    $word.UserName = "Loni"
    $newText = " I choose my CUCUMBER MUHAHA" + $quote
    $rng.InsertAfter($newText)

    # Move the "_GoBack" bookmark (last-edit marker) into the freshly
    # typed sentence, the way Word itself relocates it while you type.
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }
    $goBackPos = $insertStart + (" I choose my C").Length
    $goBackRange = $d.Range($goBackPos, $goBackPos)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}
